$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Sony 80 cm (32 inches) Bravia KLV-32R302E HD Ready LED TV'
$ws.Range("B1").Value = '₹23,490.₹23,490'
$ws.Range("A2").Value = 'Sony 59.9 cm (24 inches) Bravia KLV-24P413D HD Ready LED TV (Black)'
$ws.Range("B2").Value = '₹14,490.₹14,490'
$ws.Range("A3").Value = 'Sony 80 cm (32 inches) Bravia KLV-32W512D HD Ready Smart LED TV'
$ws.Range("B3").Value = '₹27,990.₹27,990'
$ws.Range("A4").Value = 'Sony 80.1 cm (32 inches) Bravia KLV-32W672E Full HD LED Smart TV (Black)'
$ws.Range("B4").Value = '₹34,990.₹34,990'
$ws.Range("A5").Value = 'Sony 108 cm (43 inches) Bravia KD-43X7500E 4K UHD LED Smart TV (Black)'
$ws.Range("B5").Value = '₹64,500.₹64,500'
$ws.Range("A6").Value = 'Sony 101.4 cm (40 inches) KLV-40W672E Full HD LED Smart TV'
$ws.Range("B6").Value = '₹47,489.₹47,489'
$ws.Range("A7").Value = 'Sony 101.6 cm (40 inches) Bravia KLV-40W562D Full HD LED Smart TV (Black)'
$ws.Range("B7").Value = '₹45,990.₹45,990'
$ws.Range("A8").Value = 'Sony 108 cm (43 inches) Bravia KD-43X7002E 4K UHD LED Smart TV'
$ws.Range("B8").Value = '₹60,500.₹60,500'
$ws.Range("A9").Value = 'Sony 123.2 cm (49 inches) Bravia KLV-49W672E Full HD Smart LED TV'
$ws.Range("B9").Value = '₹61,400.₹61,400'
$ws.Range("A10").Value = 'Sony 72.4 cm (29 inches) BRAVIA KLV-29P423D HD Ready LED TV'
$ws.Range("B10").Value = '₹19,990.₹19,990'
$ws.Range("A11").Value = 'Sony 80 cm (32 inches) Bravia KLV-32R412D HD Ready LED TV'
$ws.Range("B11").Value = '₹25,990.₹25,990'
$ws.Range("A12").Value = 'Sony 108cm (43 inches) KLV-43W772E Full HD LED Smart TV'
$ws.Range("B12").Value = '₹52,350.₹52,350'
$ws.Range("A13").Value = 'Sony 101.6 cm (40 inches) Bravia KLV-40W562D Full HD Smart LED TV'
$ws.Range("B13").Value = '₹49,999.₹49,999'
$ws.Range("A14").Value = 'Sony 80 cm (32 inches) Bravia KLV-32W622E HD Ready LED Smart TV (Black)'
$ws.Range("B14").Value = '₹30,990.₹30,990'
$ws.Range("A15").Value = 'Sony 138.8 cm (55 inches) Bravia KD-55X7002E 4K UHD LED Smart TV'
$ws.Range("B15").Value = '₹98,400.₹98,400'
$ws.Range("A16").Value = 'Sony 80.0 cm (32 inches) KLV-W512D HD Ready LED Smart TV (Black)'
$ws.Range("B16").Value = '  29,500'
$ws.Range("A17").Value = 'Sony 108 cm (43 inches) Bravia KDL-43W800D Full HD 3D LED Android TV'
$ws.Range("B17").Value = '₹63,900.₹63,900'
$ws.Range("A18").Value = 'Sony 123.2 cm (49 inches) BRAVIA KLV-49W772E Full HD Smart LED TV'
$ws.Range("B18").Value = '₹69,000.₹69,000'
$ws.Range("A19").Value = 'Sony 101.6 cm (40 inches) Bravia KLV-40R352E Full HD LED TV'
$ws.Range("B19").Value = '₹41,990.₹41,990'
$ws.Range("A20").Value = 'Sony 138.8 cm (55 inches) Bravia KD-55X8200E 4K UHD LED Smart TV'
$ws.Range("B20").Value = '₹1,12,970.₹1,12,970'

# Extend used range to include column C (matches dimension A1:C20 in target)
$ws.Cells.Item(20, 3).Font.Bold = $false

# Restore the selection that was active when the workbook was saved
[void]$ws.Range("C7").Select()
